$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Price" column (D) values; force text to avoid numeric auto-conversion
# while keeping the default (unstyled) cell format, matching the source workbook.
$priceUpdates = @{
    2 = "71.038.07"
    3 = "3.831.41"
    5 = "706.15"
    6 = "171.22"
    7 = "3.829.16"
    13 = "0.0000255"
    14 = "36.50"
    15 = "4.477.09"
    16 = "3.852.01"
    17 = "71.110.89"
    18 = "7.21"
    20 = "17.38"
    21 = "495.18"
    22 = "10.63"
    24 = "85.47"
    27 = "12.10"
    28 = "3.984.66"
    31 = "3.09"
    33 = "2.22"
    35 = "0.176"
    36 = "3.800.12"
    41 = "2.32"
    43 = "3.32"
    47 = "163.87"
    48 = "430.27"
    49 = "48.91"
    51 = "1.37"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.Style = "Normal"
}

# Update "Volume(1h)" column (E) values (plain text, keep spacing)
$volumeUpdates = @{
    2 = "  -0.02%  "
    3 = "  -0.71%  "
    4 = "  -0.08%  "
    5 = "  +1.86%  "
    6 = "  -0.98%  "
    7 = "  -0.72%  "
    8 = "  -0.04%  "
    9 = "  -0.46%  "
    10 = "  -1.16%  "
    11 = "  +0.24%  "
    12 = "  -0.64%  "
    13 = "  -1.72%  "
    14 = "  -0.45%  "
    15 = "  -0.87%  "
    16 = "  -0.41%  "
    17 = "  -0.04%  "
    18 = "  -0.53%  "
    19 = "  +0.08%  "
    20 = "  -2.53%  "
    21 = "  +1.57%  "
    22 = "  -4.99%  "
    23 = "  +1.91%  "
    24 = "  +0.87%  "
    25 = "  -2.16%  "
    26 = "  +0.75%  "
    27 = "  -2.44%  "
    28 = "  -0.73%  "
    29 = "  -3.13%  "
    30 = "  -0.11%  "
    31 = "  -0.85%  "
    32 = "  -2.80%  "
    33 = "  -3.42%  "
    34 = "  -1.43%  "
    35 = "  -2.67%  "
    36 = "  -0.30%  "
    37 = "  -1.26%  "
    38 = "  -0.04%  "
    39 = "  -1.73%  "
    40 = "  +3.88%  "
    41 = "  -2.81%  "
    42 = "  -1.23%  "
    43 = "  -3.67%  "
    45 = "  +0.02%  "
    46 = "  +1.15%  "
    47 = "  -0.25%  "
    48 = "  +4.05%  "
    49 = "  +0.46%  "
    50 = "  +0.41%  "
    51 = "  -1.56%  "
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $volumeUpdates[$row]
}
